# Auto-generated Excel COM-interop script to apply the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are forced to Text so numeric-looking strings
# ('5.32', '0.998', etc.) are not silently coerced into floating point numbers,
# which would lose the authored formatting (trailing zeros, thousand-dot grouping).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.316.25'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.488.83'
$ws.Range("D3").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '521.56'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.12'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.560'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.508.47'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0992'
$ws.Range("D10").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.32'
$ws.Range("D12").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.930.59'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.253.00'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.18'
$ws.Range("D16").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.494.80'
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.73'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '322.26'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.20'
$ws.Range("D21").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.78'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.79'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.413'
$ws.Range("D25").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("D27").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0754'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.89'
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.37'
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.19'
$ws.Range("D32").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.15'
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.32'
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.04'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.67'
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.48'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.802'
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.20'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '277.21'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.47'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.601'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.68'
$ws.Range("D46").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0494'
$ws.Range("D48").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.24'
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.745.83'
$ws.Range("D51").Style = "Normal"

# Coin name / link / volume columns (B, C, E) are never numeric-looking,
# so a plain .Value assignment keeps them as text with original styling.
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("E16").Value = '  -2.43%  '
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("E19").Value = '  -2.75%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("E26").Value = '  -1.07%  '
$ws.Range("E27").Value = '  -0.80%  '
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("E32").Value = '  +5.46%  '
$ws.Range("E33").Value = '  -2.16%  '
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  +4.68%  '
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("E44").Value = '  -2.73%  '
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  -2.96%  '
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("E48").Value = '  -0.78%  '
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("E51").Value = '  -0.17%  '
